$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Risk"
$ws.Range("C1").Value = "Cell"
$ws.Range("D1").Value = "Mask"
$ws.Range("A1").Value = "Trade ID"

$ws.Range("D10").WrapText = $true
